# Insert a new data row at row 54 (this shifts the existing rows 54-135
# down to 55-136, extending the sheet's dimension from A1:R135 to A1:R136).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

# Populate the newly inserted row with the new "Cebollín baby" record.
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 45117
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100112038
$ws.Range("G54").Value = "Cebollín baby"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 300
$ws.Range("K54").Value = 1300
$ws.Range("L54").Value = 1500
$ws.Range("M54").Value = 1367
$ws.Range("N54").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 684
$ws.Range("Q54").Value = 2
$ws.Range("R54").Value = "Hortaliza"
